$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.617.24"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.290.77"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "114.22"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +18.55%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "268.70"
$c.Style = $origStyle
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.30%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.616"
$c.Style = $origStyle
$ws.Range("E9").Value = "  +1.25%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "48.00"
$c.Style = $origStyle
$ws.Range("E10").Value = "  +5.23%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0937"
$c.Style = $origStyle
$ws.Range("E11").Value = "  +0.44%  "
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.64"
$c.Style = $origStyle
$ws.Range("E12").Value = "  +11.23%  "
$ws.Range("E13").Value = "  +1.74%  "
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "15.56"
$c.Style = $origStyle
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "2.634.56"
$ws.Range("E15").Value = "  +0.14%  "
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.847"
$c.Style = $origStyle
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "2.292.65"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "43.611.06"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +2.51%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.55"
$c.Style = $origStyle
$ws.Range("E20").Value = "  +5.55%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "72.47"
$c.Style = $origStyle
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  +2.71%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "233.34"
$c.Style = $origStyle
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +4.29%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = $origStyle
$ws.Range("E25").Value = "  +14.51%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +3.88%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "42.03"
$c.Style = $origStyle
$ws.Range("E28").Value = "  +4.83%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.41"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -0.02%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "176.23"
$c.Style = $origStyle
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.Style = $origStyle
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "21.53"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -1.11%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = $origStyle
$ws.Range("E34").Value = "  +3.69%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.126"
$c.Style = $origStyle
$ws.Range("E35").Value = "  +0.86%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +9.53%  "
$ws.Range("E37").Value = "  +0.68%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0359"
$c.Style = $origStyle
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("E39").Value = "  +13.14%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "13.95"
$c.Style = $origStyle
$ws.Range("E40").Value = "  +13.85%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "74.27"
$c.Style = $origStyle
$ws.Range("E41").Value = "  +14.97%  "
$ws.Range("E42").Value = "  +3.61%  "
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = $origStyle
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("E44").Value = "  +7.03%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = $origStyle
$ws.Range("E45").Value = "  +14.98%  "
$ws.Range("E46").Value = "  +0.20%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.75"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -0.40%  "
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "102.64"
$c.Style = $origStyle
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("E50").Value = "  +3.90%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.447"
$c.Style = $origStyle
$ws.Range("E51").Value = "  +4.70%  "
